# Creditos.xlsx - add a new data row (row 5) for a payment/transfer scenario
# from an account with no balance, mirroring the structure of existing row 4.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Duplicate row 4 (values + formatting) into the new row 5
$ws.Range("A4:T4").Copy($ws.Range("A5:T5"))

# New scenario id for this row
$ws.Range("A5").Value2 = 4

# Re-create the hyperlink on the new row's correo (N column), matching the
# other rows' mailto hyperlink
$ws.Hyperlinks.Add($ws.Range("N5"), "mailto:jalzate@todo1.net")

# Taller rows (2-5) to accommodate the data-driven scenario rows
$ws.Rows.Item(2).RowHeight = 30
$ws.Rows.Item(3).RowHeight = 30
$ws.Rows.Item(4).RowHeight = 30
$ws.Rows.Item(5).RowHeight = 30

# Update view state: zoom and active selection
$ws.Activate()
$excel.ActiveWindow.Zoom = 98
$ws.Range("A6").Select()
